$p = $ppt.ActivePresentation
$s = $p.Slides.Item(11)
$tblShape = $s.Shapes.Item(2)
$tbl = $tblShape.Table
$newRow = $tbl.Rows.Add(5)
Write-Host $tbl.Rows.Count
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $t = $tbl.Cell($i, 1).Shape.TextFrame.TextRange.Text
    Write-Host $i "len=" $t.Length
}
